# Add age-group to the Data Wrangling area (slide 7)
$p = $ppt.ActivePresentation

# --- Slide 7: "Data Wrangling" bullet gains an "Age-Group" sub-bullet ---
$s7 = $p.Slides.Item(7)
$contentShape = $s7.Shapes.Item(2)
$tr7 = $contentShape.TextFrame.TextRange

# Locate the "Data Wrangling" paragraph (last paragraph in this placeholder)
$n7 = $tr7.Paragraphs().Count
$dataWranglingPara = $tr7.Paragraphs($n7, 1)

# Split "Data Wrangling" into two runs: "Data " and "Wrangling"
$firstPart = $dataWranglingPara.Characters(1, 5)
$firstPart.Text = "Data "

# Append a new sub-bullet paragraph "Age-Group" at the next indent level
$inserted = $tr7.InsertAfter("`rAge-Group")
$n7b = $tr7.Paragraphs().Count
$ageGroupPara = $tr7.Paragraphs($n7b, 1)
$ageGroupPara.IndentLevel = 2

# --- Slide 15: merge "groups " and "variances" runs into a single run ---
$s15 = $p.Slides.Item(15)
$chiSquareShape = $s15.Shapes.Item(2)
$tr15 = $chiSquareShape.TextFrame.TextRange

$n15 = $tr15.Paragraphs().Count
$variancesPara = $tr15.Paragraphs($n15, 1)

# "Compared two or more " is 22 characters; "groups variances" is the remaining 17
$prefixLen = "Compared two or more ".Length
$mergeLen = $variancesPara.Length - $prefixLen
$mergeRange = $variancesPara.Characters($prefixLen + 1, $mergeLen)
$mergeRange.Text = "groups variances"
